$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# "Technology Stack" sub-items (1-based paragraph indices) before the edit:
#  6: Platform: Azure Virtual Desktop (AVD)
#  7: Session Hosts: Windows 10/11 Multi-session (D4s_v5 VMs)
#  8: Storage: Azure Files + FSLogix
#  9: Identity: Azure AD / Microsoft Entra ID
# 10: Security: Network Security Groups, Private Endpoints
# 11: Monitoring: Azure Monitor, Diagnostics Services
#
# Consolidate down to 3 items: Platform, Storage & Identity, Security & Monitoring.
# For each merge pair we rewrite the text of the *second* (later) paragraph and
# delete the *first* (earlier) one, so we never delete the text frame's very
# last paragraph (which would leave a stray empty trailing paragraph behind).

# Security & Monitoring (paragraphs 10 + 11) -> keep paragraph 11, drop 10.
$tr.Paragraphs(11, 1).Runs(1, 1).Text = "Security & Monitoring: Network Security Groups, Private Endpoints, Azure Monitor"
$tr.Paragraphs(10, 1).Delete()

# Storage & Identity (paragraphs 8 + 9) -> keep paragraph 9, drop 8.
$tr.Paragraphs(9, 1).Runs(1, 1).Text = "Storage & Identity: Azure Files + FSLogix, Azure AD / Microsoft Entra ID"
$tr.Paragraphs(8, 1).Delete()

# Platform (paragraphs 6 + 7) -> keep paragraph 6, drop 7.
$tr.Paragraphs(6, 1).Runs(1, 1).Text = "Platform: Azure Virtual Desktop (AVD) with Windows 10/11 Multi-session hosts"
$tr.Paragraphs(7, 1).Delete()
